# Se agrega checklist de entrega y lista de parametros
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# xlPasteFormats constant
$xlPasteFormats = -4122
$xlPasteAll = -4104
$xlCenter = -4108
$xlNone = -4142
$xlContinuous = 1
$xlThin = 2

# ---------------------------------------------------------------------------
# 1) Insert a new blank row above old row 2. This shifts:
#      old row 2            -> row 3   ("Lista de entrada y salidas del PLC")
#      old rows 4..26       -> rows 5..27 (the two parameter tables)
#    and leaves row 4 empty (gap), matching the target layout.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).Insert()

# ---------------------------------------------------------------------------
# 2) Re-apply the title-bar look (fill + centered + border) freshly to B1:F1
#    and B3:F3 so Excel records brand-new style entries for them (mirrors the
#    authors' edit, which re-saved these ranges as part of the restructure).
# ---------------------------------------------------------------------------
$ws.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial($xlPasteFormats)

$ws.Range("B3:F3").Copy()
$ws.Range("B3:F3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Build the new separator row 2 (A2:G2): plain cells on B:F with a thin
#    top+bottom border (no fill), and "border-only" cells on A2/G2.
# ---------------------------------------------------------------------------
$row2 = $ws.Range("A2:G2")
$row2.Value2 = ""
$row2.Interior.Pattern = $xlNone
$row2.Borders.Item(7).LineStyle = $xlNone   # xlInsideVertical off (left edge default)
$row2.Borders.Item(8).LineStyle = $xlNone
$bf = $ws.Range("B2:F2")
$bf.Borders.Item(9).LineStyle = $xlContinuous  # xlEdgeTop
$bf.Borders.Item(9).Weight = $xlThin
$bf.Borders.Item(10).LineStyle = $xlContinuous # xlEdgeBottom
$bf.Borders.Item(10).Weight = $xlThin
$bf.HorizontalAlignment = $xlCenter

$a2g2 = $ws.Range("A2")
$a2g2.Borders.Item(7).LineStyle = $xlNone
$g2 = $ws.Range("G2")
$g2.Borders.Item(7).LineStyle = $xlNone

# ---------------------------------------------------------------------------
# 4) Style fix-ups inside the parameter tables (these cells swap which
#    alternating fill they use):
# ---------------------------------------------------------------------------
$ws.Range("B9").Copy()
$ws.Range("C9").PasteSpecial($xlPasteFormats)
$ws.Range("C9").Value2 = "MaquinaAbajo"

$ws.Range("B11").Copy()
$ws.Range("C11").PasteSpecial($xlPasteFormats)
$ws.Range("C11").Value2 = "PuertaCerrada"

$ws.Range("F18").Copy()
$ws.Range("F18").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

Write-Output "done"
